# Add 2 new books to the list (rows 66 and 67), pushing the trailing
# blank/formatting row down to row 68.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Move the existing trailing blank formatting row (currently row 66,
#    all four cells styled but empty) down to the new last row, 68.
$ws.Range("E66:H66").Copy()
$ws.Range("E68:H68").PasteSpecial(-4122)

# 2) Stamp rows 66 and 67 with the same cell formatting used by the other
#    data rows (copy format only from row 65, the last data row) so they
#    reuse the existing styles instead of creating new ones.
$ws.Range("E65:H65").Copy()
$ws.Range("E66:H67").PasteSpecial(-4122)
$ws.Rows("66:67").RowHeight = 21

# 3) Fill in the two new book entries.
$ws.Range("E66").Value = 62
$ws.Range("F66").Value = "Solitaire"
$ws.Range("G66").Value = "Alice oslem"
$ws.Range("H66").Value = "Fiction"

$ws.Range("E67").Value = 63
$ws.Range("F67").Value = "The Silva - Mind controlling techniques"
$ws.Range("G67").Value = "Jose Silva"
$ws.Range("H67").Value = "Self-Decvelopment"

# 4) Update the active selection to reflect where the editor left off.
$null = $ws.Range("J56").Select()
